# Update "想去人数" (F column) counts across sheets as published at
# generated output 7921097 (gh-pages update).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 1203
$ws1.Range("F12").Value = 686
$ws1.Range("F23").Value = 1230
$ws1.Range("F24").Value = 295
$ws1.Range("F28").Value = 2503
$ws1.Range("F41").Value = 203

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 34
$ws2.Range("F23").Value = 421

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 2131
$ws3.Range("F11").Value = 773
$ws3.Range("F12").Value = 112

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2131
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 1203
$ws4.Range("F27").Value = 1230
$ws4.Range("F28").Value = 295
$ws4.Range("F32").Value = 2503
$ws4.Range("F49").Value = 203
